$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('E2').Value = '2026-02-17 17:18:30'
$ws.Range('K2').Value = '6.9 MJ/m2'
$ws.Range('E3').Value = '2026-02-17 17:18:32'
$ws.Range('K3').Value = '6.7 MJ/m2'
$ws.Range('O3').Value = '-4.5 °C'
$ws.Range('E4').Value = '2026-02-17 17:18:35'
$ws.Range('H4').Value = '80%'
$ws.Range('K4').Value = '7.4 MJ/m2'
$ws.Range('E5').Value = '2026-02-17 17:18:37'
$ws.Range('H5').Value = '94%'
$ws.Range('K5').Value = '5.9 MJ/m2'
$ws.Range('O5').Value = '-4.0 °C'
$ws.Range('E6').Value = '2026-02-17 17:18:39'
$ws.Range('J6').Value = '1018.1 hPa'
$ws.Range('K6').Value = '9.8 MJ/m2'
$ws.Range('O6').Value = '10.6 °C'
$ws.Range('E7').Value = '2026-02-17 17:18:42'
$ws.Range('H7').Value = '58%'
$ws.Range('K7').Value = '13.1 MJ/m2'
$ws.Range('E8').Value = '2026-02-17 17:18:44'
$ws.Range('J8').Value = '1017.9 hPa'
$ws.Range('O8').Value = '10.8 °C'
$ws.Range('E9').Value = '2026-02-17 17:18:46'
$ws.Range('H9').Value = '54%'
$ws.Range('K9').Value = '8.7 MJ/m2'
$ws.Range('E10').Value = '2026-02-17 17:18:49'
$ws.Range('K10').Value = '10.6 MJ/m2'
$ws.Range('O10').Value = '10.9 °C'
$ws.Range('E11').Value = '2026-02-17 17:18:51'
$ws.Range('O11').Value = '7.8 °C'
$ws.Range('E12').Value = '2026-02-17 17:18:53'
$ws.Range('E13').Value = '2026-02-17 17:18:55'
$ws.Range('J13').Value = '1017.5 hPa'
$ws.Range('K13').Value = '9.4 MJ/m2'
$ws.Range('O13').Value = '7.1 °C'
$ws.Range('E14').Value = '2026-02-17 17:18:58'
$ws.Range('E15').Value = '2026-02-17 17:19:00'
$ws.Range('E16').Value = '2026-02-17 17:19:02'
$ws.Range('H16').Value = '68%'
$ws.Range('E17').Value = '2026-02-17 17:19:05'
$ws.Range('H17').Value = '78%'
$ws.Range('K17').Value = '10.4 MJ/m2'
$ws.Range('E18').Value = '2026-02-17 17:19:07'
$ws.Range('J18').Value = '1018.3 hPa'
$ws.Range('K18').Value = '10.5 MJ/m2'
$ws.Range('O18').Value = '10.2 °C'
$ws.Range('E19').Value = '2026-02-17 17:19:10'
$ws.Range('K19').Value = '9.4 MJ/m2'
$ws.Range('E20').Value = '2026-02-17 17:19:12'
$ws.Range('H20').Value = '58%'
$ws.Range('K20').Value = '9.8 MJ/m2'
$ws.Range('O20').Value = '-2.1 °C'
$ws.Range('E21').Value = '2026-02-17 17:19:15'
$ws.Range('K21').Value = '6.4 MJ/m2'
$ws.Range('O21').Value = '9.8 °C'
$ws.Range('E22').Value = '2026-02-17 17:19:17'
$ws.Range('E23').Value = '2026-02-17 17:19:19'
$ws.Range('K23').Value = '12.2 MJ/m2'
$ws.Range('E24').Value = '2026-02-17 17:19:21'
$ws.Range('K24').Value = '14.9 MJ/m2'
$ws.Range('O24').Value = '12.6 °C'
$ws.Range('E25').Value = '2026-02-17 17:19:24'
$ws.Range('K25').Value = '13.5 MJ/m2'
$ws.Range('O25').Value = '-1.0 °C'
$ws.Range('E26').Value = '2026-02-17 17:19:26'
$ws.Range('E27').Value = '2026-02-17 17:19:28'
$ws.Range('H27').Value = '50%'
$ws.Range('O27').Value = '-0.6 °C'
$ws.Range('E28').Value = '2026-02-17 17:19:31'
$ws.Range('O28').Value = '8.7 °C'
$ws.Range('E29').Value = '2026-02-17 17:19:33'
$ws.Range('K29').Value = '10.0 MJ/m2'
$ws.Range('E30').Value = '2026-02-17 17:19:35'
$ws.Range('H30').Value = '59%'
$ws.Range('J30').Value = '1018.0 hPa'
$ws.Range('K30').Value = '7.9 MJ/m2'
$ws.Range('E31').Value = '2026-02-17 17:19:38'
$ws.Range('J31').Value = '1018.2 hPa'
$ws.Range('K31').Value = '9.5 MJ/m2'
$ws.Range('E32').Value = '2026-02-17 17:19:40'
$ws.Range('K32').Value = '11.1 MJ/m2'
$ws.Range('O32').Value = '8.4 °C'
$ws.Range('E33').Value = '2026-02-17 17:19:43'
$ws.Range('H33').Value = '39%'
$ws.Range('K33').Value = '7.2 MJ/m2'
$ws.Range('E34').Value = '2026-02-17 17:19:45'
$ws.Range('H34').Value = '47%'
$ws.Range('K34').Value = '11.2 MJ/m2'
$ws.Range('O34').Value = '1.2 °C'
$ws.Range('E35').Value = '2026-02-17 17:19:47'
$ws.Range('K35').Value = '9.5 MJ/m2'
$ws.Range('O35').Value = '7.2 °C'
$ws.Range('E36').Value = '2026-02-17 17:19:50'
$ws.Range('K36').Value = '10.6 MJ/m2'
$ws.Range('E37').Value = '2026-02-17 17:19:52'
$ws.Range('H37').Value = '69%'
$ws.Range('J37').Value = '1018.6 hPa'
$ws.Range('E38').Value = '2026-02-17 17:19:54'
$ws.Range('O38').Value = '11.4 °C'
$ws.Range('E39').Value = '2026-02-17 17:19:56'
$ws.Range('H39').Value = '55%'
$ws.Range('K39').Value = '10.5 MJ/m2'
$ws.Range('O39').Value = '-2.6 °C'
$ws.Range('E40').Value = '2026-02-17 17:19:59'
$ws.Range('E41').Value = '2026-02-17 17:20:01'
$ws.Range('K41').Value = '12.3 MJ/m2'
$ws.Range('L41').Value = '37.4 km/h - 298º 16:47 TU'
$ws.Range('O41').Value = '16.7 °C'
$ws.Range('E42').Value = '2026-02-17 17:20:04'
$ws.Range('E43').Value = '2026-02-17 17:20:06'
$ws.Range('K43').Value = '13.3 MJ/m2'
$ws.Range('O43').Value = '8.0 °C'
$ws.Range('E44').Value = '2026-02-17 17:20:08'
$ws.Range('K44').Value = '10.2 MJ/m2'
$ws.Range('O44').Value = '-3.3 °C'
$ws.Range('E45').Value = '2026-02-17 17:20:11'
$ws.Range('H45').Value = '65%'
$ws.Range('E46').Value = '2026-02-17 17:20:13'
$ws.Range('H46').Value = '54%'
$ws.Range('K46').Value = '14.3 MJ/m2'
$ws.Range('O46').Value = '15.7 °C'
